$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Rebuild the whole textbox content with the new paragraphs. Paragraphs are
# separated with the PowerPoint paragraph-mark character (`r).
$restApiLink  = "https://www.youtube.com/watch?v=ST8XxjOTIsg&list=PLTCrU9sGybupzS5-3iYTsYUI1emBDKdHu"
$sudocodeLink = "https://www.youtube.com/@sudocode/playlists"
$codingLink   = "https://www.youtube.com/@LeadCodingbyFRAZ"

$paragraphs = @(
    "System Design, Backend Stack , Frontend Stack",
    "",
    "For Revision of REST API, follow the below link:",
    $restApiLink,
    "",
    "",
    $sudocodeLink,
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "Coding ",
    $codingLink,
    "",
    "",
    "",
    ""
)

$tr.Text = [string]::Join([char]13, $paragraphs)

# Re-apply the hyperlinks lost when the text was rewritten, reusing the
# existing relationship for the REST API link and creating new ones for the
# two newly-added links. Offsets are tracked manually (1-based, PowerPoint
# convention) instead of relying on TextRange.Paragraphs(), which undercounts
# trailing empty paragraphs in this host.
$offset = 1
foreach ($paraText in $paragraphs) {
    if ($paraText -eq $restApiLink -or $paraText -eq $sudocodeLink -or $paraText -eq $codingLink) {
        $target = $tr.Characters($offset, $paraText.Length)
        $target.ActionSettings(1).Hyperlink.Address = $paraText
    }
    $offset += $paraText.Length + 1
}
